$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Symmetric Costs" table (rows 6-7) ---
# SVM row (6): values stay the same, but it is no longer the highlighted
# (best) result, so its fill is cleared back to "no fill".
$ws.Range("C6:E6").Interior.Pattern = -4142  # xlPatternNone

# MLP row (7): now holds the best results -> new values, fill stays as-is
# (it already carries the highlighted/shaded style).
$ws.Range("C7").Value = 0.81699999999999995
$ws.Range("D7").Value = 0.86
$ws.Range("E7").Value = 0.86

# --- Update the "Asymmetric Costs" table (rows 16-17) ---
# SVM row (16): same story - clear the highlight fill, keep the values.
$ws.Range("C16:E16").Interior.Pattern = -4142  # xlPatternNone

# MLP row (17): new best values.
$ws.Range("C17").Value = 0.80300000000000005
$ws.Range("D17").Value = 0.86
$ws.Range("E17").Value = 0.86

# --- Misc: move the active selection cursor like in the saved file ---
$ws.Range("I15").Select()
